$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone header style (bold/border/center-top alignment) from H1 onto I1/J1
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-28: I column = 1, J column = copy of H column value
for ($r = 2; $r -le 28; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
